$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [string][char]0x2705 + " 1000 Bs = 7.35 = 30095.59 pesos"
$newLine1 = [string][char]0x2705 + " 1000 Bs = 7.46 = 30597.01 pesos"
$oldLine2 = [string][char]0x2705 + " 30095.59 pesos = 7.31 = 961.11 Bs"
$newLine2 = [string][char]0x2705 + " 30597.01 pesos = 7.45 = 965.24 Bs"

$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$cellA1.Value = $text

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 134
$wsTasas.Range("O10").Value = 4100
$wsTasas.Range("N12").Value = 4105
$wsTasas.Range("O12").Value = 129.5
